$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Wipe all cell VALUES (keeps the existing styles / shared-string
#    table is reset so we can rebuild it from scratch in a very
#    specific order - this lets us reproduce the exact sharedStrings
#    ordering required by the target workbook).
# ------------------------------------------------------------------
$ws.Cells.ClearContents()

# ------------------------------------------------------------------
# 2) Give the 3 brand-new index cells (A17:A19) the same style as the
#    rest of column A (bold / bordered / centered) by copying the
#    formatting from A16, which already carries that style.
# ------------------------------------------------------------------
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3) Re-enter every value. The ORDER in which new text values are
#    assigned controls the order they are appended to the shared
#    string table, so we deliberately go column B top-to-bottom
#    first (this introduces HKL, the scheme names - including the
#    newly-run Gaussian-Quadrature / Spiral schemes - in the exact
#    order needed), and only afterwards fill the C2:M2 header labels.
# ------------------------------------------------------------------

# --- Column B (row labels), rows 2-19, in fill order -------------
$ws.Range("B2").Value  = "HKL"
$ws.Range("B3").Value  = "ND Single"
$ws.Range("B4").Value  = "RD Single"
$ws.Range("B5").Value  = "TD Single"
$ws.Range("B6").Value  = "Morris"
$ws.Range("B7").Value  = "Ring Perpendicular to ND"
$ws.Range("B8").Value  = "Ring Perpendicular to RD"
$ws.Range("B9").Value  = "Ring Perpendicular to TD"
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

# --- Row 2 data-column headers (HKL family labels), C2:M2 ---------
$ws.Range("C2").Value = "[1, 1, 0]"
$ws.Range("D2").Value = "[2, 0, 0]"
$ws.Range("E2").Value = "[2, 1, 1]"
$ws.Range("F2").Value = "[2, 2, 0]"
$ws.Range("G2").Value = "[3, 1, 0]"
$ws.Range("H2").Value = "[2, 2, 2]"
$ws.Range("I2").Value = "[3, 2, 1]"
$ws.Range("J2").Value = "[4, 0, 0]"
$ws.Range("K2").Value = "2Pairs"
$ws.Range("L2").Value = "4Pairs"
$ws.Range("M2").Value = "MaxUnique"

# --- Row 1 numeric column headers (0-11) --------------------------
for ($c = 2; $c -le 13; $c++) {
    $ws.Cells.Item(1, $c).Value = $c - 2
}

# --- Column A numeric row indices (0-17) for rows 2-19 ------------
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# --- Data body: C3:M16 (existing rows) and C17:M19 (new rows) -----
for ($r = 3; $r -le 19; $r++) {
    for ($c = 3; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}
